# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E21) listed the mora periods out of
# order (1701, 1612, 1611, 1610, 1609, 1608). Re-sort them ascending
# (1608 .. 1701) to reflect the refreshed account-statement database.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periodos = @("1608", "1609", "1610", "1611", "1612", "1701")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
}
